# Update the "想去人数" (interested count) values in column F for the
# data sheets "展览" and "全部类型". Both sheets contain the same
# table contents and both need the identical updates.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 367
    3  = 360
    4  = 1870
    11 = 4456
    14 = 1241
    16 = 51
    17 = 815
    19 = 440
    21 = 217
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
